$d = $word.ActiveDocument

# The template has four repeated panels. Panel 1's "Part #1"/"TAG1" runs
# are already merged into single runs. Panels 2-4 still have the label
# ("Part #" / "TAG") and the trailing panel number split across two
# separate <w:r> runs with identical run properties. Re-running Find &
# Replace over the full (concatenated) text merges the matched span back
# into a single run, which is exactly the consolidation the diff shows.

$targets = @("Part #2", "TAG2", "Part #3", "TAG3", "Part #4", "TAG4")

foreach ($t in $targets) {
    $rng = $d.Content
    $rng.Find.Execute($t, $true, $false, $false, $false, $false, $true, 1, $false, $t, 2) | Out-Null
}

$d.Saved = $false
